$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tiny floating-point precision refresh on a handful of existing
#     timestamp cells (column A) picked up when the sheet was re-saved. ---
$ws.Range("A719").Value2 = 45914.49244710648
$ws.Range("A720").Value2 = 45914.49260667824
$ws.Range("A740").Value2 = 45914.49333403935
$ws.Range("A759").Value2 = 45914.49444172454
$ws.Range("A774").Value2 = 45914.49581134259
$ws.Range("A793").Value2 = 45914.49678841436

# --- New event rows appended after row 811 (sound alarm for overload). ---

function Set-EventRow($r, $serial, $b, $c, $d, $e) {
    $ws.Cells.Item($r, 1).Value2 = $serial
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
}

function Set-RampRow($r, $timestamp, $b) {
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = "ramp"
}

Set-EventRow 812 45914.5060709375  264 660 "12:08" "Increase Unit load to 660.0 MW/Tăng tải lên 660.0 MW"
Set-EventRow 813 45914.50642701389 264 660 "00:00" "Increase Unit load to 660.0 MW/Tăng tải lên 660.0 MW"
Set-EventRow 814 45914.50908366948 264 660 "23:12" "Increase Unit load to 660.0 MW/Tăng tải lên 660.0 MW"
Set-EventRow 815 45914.50930866507 264 600 "12:13" "Increase Unit load to 600.0 MW/Tăng tải lên 600.0 MW"

Set-RampRow 816 "2025-09-14 12:46:00" 429
Set-RampRow 817 "2025-09-14 12:47:00" 442.2
Set-RampRow 818 "2025-09-14 12:48:00" 455.4
Set-RampRow 819 "2025-09-14 12:49:00" 468.6
Set-RampRow 820 "2025-09-14 12:50:00" 481.8
Set-RampRow 821 "2025-09-14 12:51:00" 494.9999999999999
Set-RampRow 822 "2025-09-14 12:52:00" 508.1999999999999
Set-RampRow 823 "2025-09-14 12:53:00" 521.4
Set-RampRow 824 "2025-09-14 12:54:00" 534.6
Set-RampRow 825 "2025-09-14 12:55:00" 547.8000000000001
Set-RampRow 826 "2025-09-14 12:56:00" 561.0000000000001
Set-RampRow 827 "2025-09-14 12:57:00" 574.2000000000002
Set-RampRow 828 "2025-09-14 12:58:00" 587.4000000000002
Set-RampRow 829 "2025-09-14 12:59:00" 600.6000000000003
Set-RampRow 830 "2025-09-14 13:00:00" 613.8000000000003
Set-RampRow 831 "2025-09-14 13:01:00" 627.0000000000003
Set-RampRow 832 "2025-09-14 13:02:00" 640.2000000000004
Set-RampRow 833 "2025-09-14 13:03:00" 653.4000000000004
Set-RampRow 834 "2025-09-14 13:04:00" 660

Set-EventRow 835 45914.50971411171 264 600 "23:13" "Increase Unit load to 600.0 MW/Tăng tải lên 600.0 MW"

Set-RampRow 836 "2025-09-14 23:46:00" 429
Set-RampRow 837 "2025-09-14 23:47:00" 442.2
Set-RampRow 838 "2025-09-14 23:48:00" 455.4
Set-RampRow 839 "2025-09-14 23:49:00" 468.6
Set-RampRow 840 "2025-09-14 23:50:00" 481.8
Set-RampRow 841 "2025-09-14 23:51:00" 494.9999999999999
Set-RampRow 842 "2025-09-14 23:52:00" 508.1999999999999
Set-RampRow 843 "2025-09-14 23:53:00" 521.4
Set-RampRow 844 "2025-09-14 23:54:00" 534.6
Set-RampRow 845 "2025-09-14 23:55:00" 547.8000000000001
Set-RampRow 846 "2025-09-14 23:56:00" 561.0000000000001
Set-RampRow 847 "2025-09-14 23:57:00" 574.2000000000002
Set-RampRow 848 "2025-09-14 23:58:00" 587.4000000000002
Set-RampRow 849 "2025-09-14 23:59:00" 600.6000000000003
Set-RampRow 850 "2025-09-15 00:00:00" 613.8000000000003
Set-RampRow 851 "2025-09-15 00:01:00" 627.0000000000003
Set-RampRow 852 "2025-09-15 00:02:00" 640.2000000000004
Set-RampRow 853 "2025-09-15 00:03:00" 653.4000000000004
Set-RampRow 854 "2025-09-15 00:04:00" 660
